$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the header out of column A into column B, then wipe out the
# "Impact category UUID" data column (A2:A18) -- the UUID values are no
# longer needed now that the database integration is in place.
$ws.Range("B1").Value = $ws.Range("A1").Value2
$ws.Range("A1:A18").ClearContents()

# Give the now-leftmost data column (Impact category) a bit more room.
$ws.Columns.Item(2).ColumnWidth = 13.7

# Update selection to match the recorded post-edit state.
$ws.Range("O18").Select()
